$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 155, shifting existing rows 155-158 down to 156-159.
$ws.Rows.Item(155).Insert()

# Populate the newly-inserted row 155 with the new weekly data point.
$ws.Cells.Item(155, 1).Value = 6
$ws.Cells.Item(155, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(155, 3).Value = "Metropolitana"
$ws.Cells.Item(155, 4).Value = 44595
$ws.Cells.Item(155, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(155, 5).Value = 13
$ws.Cells.Item(155, 6).Value = "Fruta"
$ws.Cells.Item(155, 7).Value = 100101
$ws.Cells.Item(155, 8).Value = "Berries"
$ws.Cells.Item(155, 9).Value = 100101004
$ws.Cells.Item(155, 10).Value = "Frambuesa"
$ws.Cells.Item(155, 11).Value = "Sin especificar"
$ws.Cells.Item(155, 12).Value = "Primera"
$ws.Cells.Item(155, 13).Value = 120
$ws.Cells.Item(155, 14).Value = 8000
$ws.Cells.Item(155, 15).Value = 8000
$ws.Cells.Item(155, 16).Value = 8000
$ws.Cells.Item(155, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(155, 18).Value = "Provincia de Linares"
$ws.Cells.Item(155, 19).Value = 4000
$ws.Cells.Item(155, 20).Value = 2
